# Fruta / hortaliza, semanal
# Insert a new weekly price record (2022-02-18) for "Terminal La Palmera de
# La Serena" / Platano, ahead of the existing rows, shifting the rest of
# the table down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above row 515 (existing data shifts down to 518+)
$ws.Rows.Item(515).Resize(3).Insert()

$newDate = "2022-02-18"

# Row 515: Pinton
$ws.Cells.Item(515, 1).Value = 8
$ws.Cells.Item(515, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(515, 3).Value = "Coquimbo"
$ws.Cells.Item(515, 4).Value = $newDate
$ws.Cells.Item(515, 5).Value = 4
$ws.Cells.Item(515, 6).Value = "Fruta"
$ws.Cells.Item(515, 7).Value = 100108
$ws.Cells.Item(515, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(515, 9).Value = 100108006
$ws.Cells.Item(515, 10).Value = "Plátano"
$ws.Cells.Item(515, 11).Value = "Sin especificar"
$ws.Cells.Item(515, 12).Value = "Pintón"
$ws.Cells.Item(515, 13).Value = 80
$ws.Cells.Item(515, 14).Value = 15000
$ws.Cells.Item(515, 15).Value = 15000
$ws.Cells.Item(515, 16).Value = 15000
$ws.Cells.Item(515, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(515, 18).Value = "Ecuador"
$ws.Cells.Item(515, 19).Value = 750
$ws.Cells.Item(515, 20).Value = 20

# Row 516: Primera Maduro
$ws.Cells.Item(516, 1).Value = 8
$ws.Cells.Item(516, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(516, 3).Value = "Coquimbo"
$ws.Cells.Item(516, 4).Value = $newDate
$ws.Cells.Item(516, 5).Value = 4
$ws.Cells.Item(516, 6).Value = "Fruta"
$ws.Cells.Item(516, 7).Value = 100108
$ws.Cells.Item(516, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(516, 9).Value = 100108006
$ws.Cells.Item(516, 10).Value = "Plátano"
$ws.Cells.Item(516, 11).Value = "Sin especificar"
$ws.Cells.Item(516, 12).Value = "Primera Maduro"
$ws.Cells.Item(516, 13).Value = 120
$ws.Cells.Item(516, 14).Value = 17000
$ws.Cells.Item(516, 15).Value = 17000
$ws.Cells.Item(516, 16).Value = 17000
$ws.Cells.Item(516, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(516, 18).Value = "Ecuador"
$ws.Cells.Item(516, 19).Value = 850
$ws.Cells.Item(516, 20).Value = 20

# Row 517: Primera Pinton
$ws.Cells.Item(517, 1).Value = 8
$ws.Cells.Item(517, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(517, 3).Value = "Coquimbo"
$ws.Cells.Item(517, 4).Value = $newDate
$ws.Cells.Item(517, 5).Value = 4
$ws.Cells.Item(517, 6).Value = "Fruta"
$ws.Cells.Item(517, 7).Value = 100108
$ws.Cells.Item(517, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(517, 9).Value = 100108006
$ws.Cells.Item(517, 10).Value = "Plátano"
$ws.Cells.Item(517, 11).Value = "Sin especificar"
$ws.Cells.Item(517, 12).Value = "Primera Pintón"
$ws.Cells.Item(517, 13).Value = 120
$ws.Cells.Item(517, 14).Value = 18000
$ws.Cells.Item(517, 15).Value = 18000
$ws.Cells.Item(517, 16).Value = 18000
$ws.Cells.Item(517, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(517, 18).Value = "Ecuador"
$ws.Cells.Item(517, 19).Value = 900
$ws.Cells.Item(517, 20).Value = 20
